$wb = $excel.ActiveWorkbook

# --- "Second" sheet: move the vertically-merged "Vertical Merge" cell from
# column C to column B (rows 10:11), and put a new "Meaningful data" value
# into the now free C11 cell. ---
$ws2 = $wb.Worksheets.Item("Second")
$ws2.Range("C10:C11").UnMerge()
$ws2.Range("C10").Cut($ws2.Range("B10"))
$ws2.Range("C11").Cut($ws2.Range("B11"))
$ws2.Range("C10").Clear()
$ws2.Range("C11").Clear()
$ws2.Range("B10:B11").Merge()
$ws2.Range("C11").Value = "Meaningful data"
$ws2.Range("C11").Select()

# --- "Main" sheet: append a new last row with a "Last cell" marker. ---
$ws1 = $wb.Worksheets.Item("Main")
$ws1.Range("B19").Value = "Last cell"

# --- Make "Main" the active sheet/tab, with the selection left on C20. ---
$ws1.Activate()
$ws1.Range("C20").Select()
